# issue #5: property land done
# Normalizes a handful of stray spaces / punctuation marks in several
# sheets, and brings the "土地" (land) sheet up to the same structured,
# machine-readable layout already used by the "股票" (stock) sheet -
# English column headers plus the property_category / category / date /
# legislator_name / legislator_id / source_file / index metadata columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 土地 (land) sheet
# ---------------------------------------------------------------
$land = $wb.Worksheets.Item("土地")

# Rename the headers to the machine-readable english names used
# elsewhere in the workbook.
$land.Range("B1").Value = "name"
$land.Range("C1").Value = "area"
$land.Range("D1").Value = "share_portion"
$land.Range("E1").Value = "owner"
$land.Range("F1").Value = "register_date"
$land.Range("G1").Value = "register_reason"
$land.Range("H1").Value = "acquire_value"

# Add the new metadata header columns (I..O).
$land.Range("I1").Value = "property_category"
$land.Range("J1").Value = "category"
$land.Range("K1").Value = "date"
$land.Range("L1").Value = "legislator_name"
$land.Range("M1").Value = "legislator_id"
$land.Range("N1").Value = "source_file"
$land.Range("O1").Value = "index"

# Copy the existing header formatting (bold / border / centered) onto
# the newly added header cells.
$land.Range("B1:H1").Copy()
$land.Range("I1").PasteSpecial(-4122)

# Clean up the stray "0612-0016 " style formatting in the land parcel
# names, and the extra space in the register dates.
$land.Range("B2").Value = "彰化縣員林鎮三條段06120016地號"
$land.Range("F2").Value = "93年03月22日"

$land.Range("B3").Value = "彰化縣員林鎮三條段06130006地號"
$land.Range("F3").Value = "93年03月22日"

$land.Range("B4").Value = "彰化縣員林鎮三條段06130012地號"
$land.Range("F4").Value = "93年03月22日"

# Fill in the new metadata columns for each data row. The "date" column
# (K) holds an ISO-looking "2012-04-30" string; mark the cell as Text
# first so the COM layer doesn't silently reinterpret it as a real date
# serial number - the formatting gets overwritten again below once the
# original row style is pasted back on top.
$land.Range("I2").Value = "land"
$land.Range("J2").Value = "normal"
$land.Range("K2").NumberFormat = "@"
$land.Range("K2").Value = "2012-04-30"
$land.Range("L2").Value = "魏明谷"
$land.Range("M2").Value = 980
$land.Range("N2").Value = "tmp386d1"
$land.Range("O2").Value = 14

$land.Range("I3").Value = "land"
$land.Range("J3").Value = "normal"
$land.Range("K3").NumberFormat = "@"
$land.Range("K3").Value = "2012-04-30"
$land.Range("L3").Value = "魏明谷"
$land.Range("M3").Value = 980
$land.Range("N3").Value = "tmp386d1"
$land.Range("O3").Value = 15

$land.Range("I4").Value = "land"
$land.Range("J4").Value = "normal"
$land.Range("K4").NumberFormat = "@"
$land.Range("K4").Value = "2012-04-30"
$land.Range("L4").Value = "魏明谷"
$land.Range("M4").Value = 980
$land.Range("N4").Value = "tmp386d1"
$land.Range("O4").Value = 16

# Copy the existing data-row formatting onto the newly added data cells.
$land.Range("B2:H2").Copy()
$land.Range("I2").PasteSpecial(-4122)

$land.Range("B3:H3").Copy()
$land.Range("I3").PasteSpecial(-4122)

$land.Range("B4:H4").Copy()
$land.Range("I4").PasteSpecial(-4122)

# ---------------------------------------------------------------
# 建物 (building) sheet
# ---------------------------------------------------------------
$building = $wb.Worksheets.Item("建物")
$building.Range("F2").Value = "93年02月17日"

# ---------------------------------------------------------------
# 汽車 (car) sheet
# ---------------------------------------------------------------
$car = $wb.Worksheets.Item("汽車")
$car.Range("B2").Value = "N16ESSENTRA"
$car.Range("E2").Value = "94年03月04日"

# ---------------------------------------------------------------
# 其他有價證券 (other valuable securities) sheet
# ---------------------------------------------------------------
$securities = $wb.Worksheets.Item("其他有價證券")
$securities.Range("B2").Value = "(九）珠寶古董字畫及#"
$securities.Range("C2").Value = "他具有相當價值之財產（總作"
$securities.Range("E2").Value = "!額：新臺幣元）"
$securities.Range("B3").Value = "財產種類"
$securities.Range("C3").Value = "項"

# ---------------------------------------------------------------
# 保險 (insurance) sheet
# ---------------------------------------------------------------
$insurance = $wb.Worksheets.Item("保險")
$insurance.Range("C2").Value = "豐碩人生終生B型"
